$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.982.73'
$ws.Range("E2").Value = '  -4.12%  '

$ws.Range("D3").Value = '2.327.82'
$ws.Range("E3").Value = '  -5.83%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.88'
$ws.Range("E5").Value = '  -4.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.63'
$ws.Range("E6").Value = '  -8.17%  '

$ws.Range("E7").Value = '  -3.80%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -5.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0813'
$ws.Range("E10").Value = '  -4.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.97'
$ws.Range("E11").Value = '  -8.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.109'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").Value = '2.692.93'
$ws.Range("E13").Value = '  -5.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.39'
$ws.Range("E14").Value = '  -6.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.64'
$ws.Range("E15").Value = '  -5.34%  '

$ws.Range("D16").Value = '2.317.23'
$ws.Range("E16").Value = '  -6.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.749'
$ws.Range("E17").Value = '  -5.09%  '

$ws.Range("D18").Value = '39.981.12'
$ws.Range("E18").Value = '  -3.96%  '

$ws.Range("D19").Value = '0.0₃0900'
$ws.Range("E19").Value = '  -4.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.06'
$ws.Range("E20").Value = '  -5.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.49'
$ws.Range("E21").Value = '  -5.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.60'
$ws.Range("E22").Value = '  -5.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.02'
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("E24").Value = '  -7.13%  '

$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  -6.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.30'
$ws.Range("E27").Value = '  -6.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -1.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("E29").Value = '  -5.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.21'
$ws.Range("E30").Value = '  -3.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '151.84'
$ws.Range("E31").Value = '  -2.29%  '

$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.09'
$ws.Range("E33").Value = '  -6.15%  '

$ws.Range("E34").Value = '  -4.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0721'
$ws.Range("E35").Value = '  -5.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  -2.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0996'
$ws.Range("E37").Value = '  -3.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.74'
$ws.Range("E38").Value = '  -5.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.69'
$ws.Range("E39").Value = '  -8.23%  '

$ws.Range("E40").Value = '  -7.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.80'
$ws.Range("E41").Value = '  -4.59%  '

$ws.Range("E42").Value = '  -6.25%  '

$ws.Range("D43").Value = '1.936.21'
$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0266'
$ws.Range("E44").Value = '  -5.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.56'
$ws.Range("E45").Value = '  -5.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.25'
$ws.Range("E46").Value = '  -1.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.67'
$ws.Range("E47").Value = '  -9.48%  '

$ws.Range("D48").Value = '2.557.51'
$ws.Range("E48").Value = '  -6.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '92.64'
$ws.Range("E49").Value = '  -4.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.89'
$ws.Range("E50").Value = '  -6.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.35'
$ws.Range("E51").Value = '  -3.20%  '
